$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
# Row 74
$ws.Range("H74").Value = 9004
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
# Row 77
$ws.Range("H77").Value = 9004
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
# Row 132
$ws.Range("H132").Value = 574.4091
$ws.Range("I132").Value = 589.8570999999999
$ws.Range("K132").Value = 1769.5713
$ws.Range("M132").Value = 760.4287000000002
# Row 137
$ws.Range("H137").Value = 1306.8572
$ws.Range("I137").Value = 854.5
$ws.Range("K137").Value = 2563.5
$ws.Range("M137").Value = -13.5
# Row 138
$ws.Range("H138").Value = 2693.0588
$ws.Range("I138").Value = 1514.2
$ws.Range("J138").Value = 3184.25
$ws.Range("K138").Value = 4542.6
$ws.Range("L138").Value = 9552.75
$ws.Range("M138").Value = 597.3999999999996
$ws.Range("N138").Value = -19832.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 14
$ws.Range("H14").Value = 1000000
$ws.Range("I14").Value = 1000000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1000000
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -999825
$ws.Range("N14").ClearContents()
# Row 16
$ws.Range("H16").Value = 25000500
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 50000000
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 50000000
$ws.Range("M16").Value = -713
$ws.Range("N16").Value = -50000574
# Row 38
$ws.Range("H38").Value = 2512500
$ws.Range("I38").Value = 25000
$ws.Range("K38").Value = 25000
$ws.Range("M38").Value = -24533
# Row 61
$ws.Range("H61").Value = 2200
$ws.Range("I61").Value = 1935.3334
$ws.Range("J61").Value = 2994
$ws.Range("K61").Value = 1935.3334
$ws.Range("L61").Value = 2994
$ws.Range("M61").Value = -1723.3334
$ws.Range("N61").Value = -3418
# Row 132
$ws.Range("H132").Value = 1420.1818
$ws.Range("I132").Value = 1539.8889
$ws.Range("K132").Value = 4619.6667
$ws.Range("M132").Value = -2089.6667
# Row 136
$ws.Range("H136").Value = 2200
$ws.Range("I136").Value = 1935.3334
$ws.Range("J136").Value = 2994
$ws.Range("K136").Value = 5806.0002
$ws.Range("L136").Value = 8982
$ws.Range("M136").Value = -3256.0002
$ws.Range("N136").Value = -14082

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 38
$ws.Range("H38").Value = 31000
$ws.Range("J38").Value = 31000
$ws.Range("L38").Value = 31000
$ws.Range("N38").Value = -31832
# Row 46
$ws.Range("H46").Value = 28124.5
$ws.Range("I46").Value = 19999
$ws.Range("J46").Value = 30833
$ws.Range("K46").Value = 19999
$ws.Range("L46").Value = 30833
$ws.Range("M46").Value = -19701
$ws.Range("N46").Value = -31429
# Row 82
$ws.Range("H82").Value = 27233.428
$ws.Range("I82").Value = 15289.083
$ws.Range("K82").Value = 15289.083
$ws.Range("M82").Value = -14906.083
# Row 85
$ws.Range("H85").Value = 27233.428
$ws.Range("I85").Value = 15289.083
$ws.Range("K85").Value = 15289.083
$ws.Range("M85").Value = -13963.083
# Row 94
$ws.Range("H94").Value = 74677.60000000001
$ws.Range("I94").Value = 92830.336
$ws.Range("J94").Value = 2066.6667
$ws.Range("K94").Value = 92830.336
$ws.Range("L94").Value = 2066.6667
$ws.Range("M94").Value = -92379.336
$ws.Range("N94").Value = -2968.6667
# Row 140
$ws.Range("H140").Value = 52093
$ws.Range("J140").Value = 53640
$ws.Range("L140").Value = 53640
$ws.Range("N140").Value = -64000

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 9
$ws.Range("H9").Value = 239499.5
$ws.Range("J9").Value = 239499.5
$ws.Range("L9").Value = 239499.5
$ws.Range("N9").Value = -239835.5
# Row 14
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
# Row 31
$ws.Range("H31").Value = 3975
$ws.Range("I31").Value = 3975
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 3975
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -3680
$ws.Range("N31").ClearContents()
# Row 34
$ws.Range("H34").Value = 3975
$ws.Range("I34").Value = 3975
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 3975
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -3773
$ws.Range("N34").ClearContents()
# Row 38
$ws.Range("H38").Value = 22333
$ws.Range("J38").Value = 21249.5
$ws.Range("L38").Value = 21249.5
$ws.Range("N38").Value = -22003.5
# Row 46
$ws.Range("H46").Value = 22333
$ws.Range("J46").Value = 21249.5
$ws.Range("L46").Value = 21249.5
$ws.Range("N46").Value = -21671.5
# Row 59
$ws.Range("H59").Value = 45000
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
# Row 62
$ws.Range("H62").Value = 3907.8333
$ws.Range("I62").Value = 3385
$ws.Range("J62").Value = 4953.5
$ws.Range("K62").Value = 3385
$ws.Range("L62").Value = 4953.5
$ws.Range("M62").Value = -2761
$ws.Range("N62").Value = -6201.5
# Row 64
$ws.Range("H64").Value = 55000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 55000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 55000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -55496
# Row 65
$ws.Range("H65").Value = 3907.8333
$ws.Range("I65").Value = 3385
$ws.Range("J65").Value = 4953.5
$ws.Range("K65").Value = 16925
$ws.Range("L65").Value = 24767.5
$ws.Range("M65").Value = -13805
$ws.Range("N65").Value = -31007.5
# Row 67
$ws.Range("H67").Value = 55000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 55000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 55000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -56716
# Row 132
$ws.Range("H132").Value = 2318.4285
$ws.Range("I132").Value = 2200.8
$ws.Range("J132").Value = 2612.5
$ws.Range("K132").Value = 6602.400000000001
$ws.Range("L132").Value = 7837.5
$ws.Range("M132").Value = -4072.400000000001
$ws.Range("N132").Value = -12897.5
# Row 134
$ws.Range("H134").Value = 3197.8
$ws.Range("I134").Value = 3500
$ws.Range("J134").Value = 2996.3333
$ws.Range("K134").Value = 10500
$ws.Range("L134").Value = 8988.999899999999
$ws.Range("M134").Value = -7965
$ws.Range("N134").Value = -14058.9999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 129
$ws.Range("H129").Value = 539.4
$ws.Range("I129").Value = 539.4
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 1618.2
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 3381.8
$ws.Range("N129").ClearContents()
# Row 133
$ws.Range("H133").Value = 2499.5
$ws.Range("I133").Value = 2499.5
$ws.Range("K133").Value = 7498.5
$ws.Range("M133").Value = -2438.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 10162.314
$ws.Range("I7").Value = 10163.482
$ws.Range("K7").Value = 10163.482
$ws.Range("M7").Value = -10051.482
# Row 93
$ws.Range("H93").Value = 66670080
$ws.Range("I93").Value = 111114616
$ws.Range("J93").Value = 3277.5
$ws.Range("K93").Value = 111114616
$ws.Range("L93").Value = 3277.5
$ws.Range("M93").Value = -111113368
$ws.Range("N93").Value = -5773.5
# Row 106
$ws.Range("H106").Value = 10924.429
$ws.Range("J106").Value = 10924.429
$ws.Range("L106").Value = 10924.429
$ws.Range("N106").Value = -13448.429
# Row 126
$ws.Range("H126").Value = 10162.314
$ws.Range("I126").Value = 10163.482
$ws.Range("K126").Value = 30490.446
$ws.Range("M126").Value = -28020.446

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1384.9
$ws.Range("I81").Value = 1384.9
$ws.Range("K81").Value = 2769.8
$ws.Range("M81").Value = -1708.8
# Row 84
$ws.Range("H84").Value = 1384.9
$ws.Range("I84").Value = 1384.9
$ws.Range("K84").Value = 13849
$ws.Range("M84").Value = -8545
# Row 126
$ws.Range("H126").Value = 4143.4546
$ws.Range("I126").Value = 3308.6667
$ws.Range("K126").Value = 9926.000100000001
$ws.Range("M126").Value = -7456.000100000001
